$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4366.788429377427
$ws.Range("C3").Value = 4366.788429377427
$ws.Range("C4").Value = 4366.788429377427
$ws.Range("C5").Value = 4218.47708522146
$ws.Range("C6").Value = 4218.47708522146
$ws.Range("C7").Value = 4218.47708522146
$ws.Range("C8").Value = 4021.217076499784
$ws.Range("C9").Value = 4021.217076499784
$ws.Range("C10").Value = 4021.217076499784
$ws.Range("C11").Value = 3948.046647506439
$ws.Range("C12").Value = 3933.02389931131
